$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Posted" header in column C (new third header next to
# Title/Article), re-using the exact same header formatting that A1/B1
# already have (bold font, thin border, centered/top alignment).
$ws.Range("C1").Value = "Posted"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# Clear the clipboard marching-ants state left by Copy().
$excel.CutCopyMode = $false

# Match the author's final selection/active-cell on the sheet.
$ws.Range("F2").Select()
